$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.140.93"
$ws.Range("E2").Value = "  -2.61%  "

# Row 3
$ws.Range("D3").Value = "1.849.59"
$ws.Range("E3").Value = "  -1.43%  "

# Row 4
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "'0.6948"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.93%  "

# Row 6
$ws.Range("D6").Value = "'238.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.75%  "

# Row 7
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
$ws.Range("D8").Value = "'0.3064"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.61%  "

# Row 9
$ws.Range("D9").Value = "'0.07542"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.35%  "

# Row 10
$ws.Range("D10").Value = "'23.47"
$ws.Range("D10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.08112"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.78%  "

# Row 12
$ws.Range("D12").Value = "1.889.47"
$ws.Range("E12").Value = "  +0.48%  "

# Row 13
$ws.Range("D13").Value = "'0.7234"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.43%  "

# Row 14
$ws.Range("D14").Value = "'5.184"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.81%  "

# Row 15
$ws.Range("D15").Value = "'89.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.19%  "

# Row 16
$ws.Range("D16").Value = "29.399.66"
$ws.Range("E16").Value = "  -1.74%  "

# Row 17
$ws.Range("D17").Value = "'5.796"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.88%  "

# Row 18
$ws.Range("D18").Value = "'241.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.91%  "

# Row 19
$ws.Range("D19").Value = "'0.000007727"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.36%  "

# Row 20
$ws.Range("D20").Value = "'13.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.36%  "

# Row 21
$ws.Range("D21").Value = "'1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.10%  "

# Row 22
$ws.Range("D22").Value = "2.137.15"
$ws.Range("E22").Value = "  -0.10%  "

# Row 23
$ws.Range("E23").Value = "  -0.10%  "

# Row 24
$ws.Range("D24").Value = "'7.641"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.58%  "

# Row 25
$ws.Range("D25").Value = "'9.031"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.85%  "

# Row 26
$ws.Range("D26").Value = "'161.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.26%  "

# Row 27
$ws.Range("D27").Value = "'0.1463"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.78%  "

# Row 28
$ws.Range("D28").Value = "'18.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.96%  "

# Row 29
$ws.Range("D29").Value = "'1.941"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.99%  "

# Row 30
$ws.Range("D30").Value = "'1.394"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.65%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.426"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.91%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.498"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.57%  "

# Row 33
$ws.Range("D33").Value = "'4.051"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.73%  "

# Row 34
$ws.Range("D34").Value = "'0.05240"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.63%  "

# Row 35
$ws.Range("D35").Value = "'1.190"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.54%  "

# Row 36
$ws.Range("D36").Value = "'0.7097"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.25%  "

# Row 37
$ws.Range("D37").Value = "'1.000"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.17%  "

# Row 38
$ws.Range("D38").Value = "'2.661"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.40%  "

# Row 39
$ws.Range("D39").Value = "'0.01862"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.98%  "

# Row 40
$ws.Range("D40").Value = "'2.698"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.07%  "

# Row 41
$ws.Range("D41").Value = "'0.9236"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.84%  "

# Row 42
$ws.Range("D42").Value = "'5.945"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.12%  "

# Row 43
$ws.Range("D43").Value = "'0.4284"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.33%  "

# Row 44
$ws.Range("D44").Value = "'70.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.19%  "

# Row 45
$ws.Range("D45").Value = "1.046.51"
$ws.Range("E45").Value = "  -6.05%  "

# Row 46
$ws.Range("D46").Value = "'1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "

# Row 47
$ws.Range("D47").Value = "'102.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.61%  "

# Row 48
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "2.031.24"
$ws.Range("E48").Value = "  -0.32%  "

# Row 49
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "'7.229"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.64%  "

# Row 50
$ws.Range("D50").Value = "'1.746"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.12%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'9.271"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.64%  "
